$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 557.44446
$ws.Range("I11").Value = 557.44446
$ws.Range("K11").Value = 557.44446
$ws.Range("M11").Value = -417.44446

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").ClearContents()
$ws.Range("N26").Value = 0

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 7849.2573
$ws.Range("J70").Value = 7935.9062
$ws.Range("L70").Value = 23807.7186
$ws.Range("N70").Value = -24347.7186

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 7849.2573
$ws.Range("J73").Value = 7935.9062
$ws.Range("L73").Value = 23807.7186
$ws.Range("N73").Value = -25679.7186

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1814.0488
$ws.Range("I98").Value = 1836
$ws.Range("J98").Value = 1611
$ws.Range("K98").Value = 1836
$ws.Range("L98").Value = 1611
$ws.Range("M98").Value = -338
$ws.Range("N98").Value = -4607

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2020.7142
$ws.Range("I106").Value = 1619.4
$ws.Range("K106").Value = 1619.4
$ws.Range("M106").Value = -988.4000000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 5913.7915
$ws.Range("I116").Value = 4630.25
$ws.Range("J116").Value = 7197.3335
$ws.Range("K116").Value = 4630.25
$ws.Range("L116").Value = 7197.3335
$ws.Range("M116").Value = -1188.25
$ws.Range("N116").Value = -14081.3335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1814.0488
$ws.Range("I122").Value = 1836
$ws.Range("J122").Value = 1611
$ws.Range("K122").Value = 5508
$ws.Range("L122").Value = 4833
$ws.Range("M122").Value = -3058
$ws.Range("N122").Value = -9733

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 44189.93
$ws.Range("I137").Value = 55342.363
$ws.Range("J137").Value = 3297.6667
$ws.Range("K137").Value = 166027.089
$ws.Range("L137").Value = 9893.000100000001
$ws.Range("M137").Value = -163477.089
$ws.Range("N137").Value = -14993.0001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 117599.336
$ws.Range("I139").Value = 120000
$ws.Range("J139").Value = 117119.2
$ws.Range("K139").Value = 120000
$ws.Range("L139").Value = 117119.2
$ws.Range("M139").Value = -114860
$ws.Range("N139").Value = -127399.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 721
$ws.Range("I4").Value = 677.6
$ws.Range("J4").Value = 775.25
$ws.Range("K4").Value = 677.6
$ws.Range("L4").Value = 775.25
$ws.Range("M4").Value = -561.6
$ws.Range("N4").Value = -1007.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 10261241
$ws.Range("I45").Value = 17096826
$ws.Range("K45").Value = 17096826
$ws.Range("M45").Value = -17096449

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1263935.2
$ws.Range("I110").Value = 1263935.2
$ws.Range("K110").Value = 1263935.2
$ws.Range("M110").Value = -1261890.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 511306.9
$ws.Range("I122").Value = 2047.9656
$ws.Range("J122").Value = 1742016
$ws.Range("K122").Value = 6143.8968
$ws.Range("L122").Value = 5226048
$ws.Range("M122").Value = -3693.8968
$ws.Range("N122").Value = -5230948

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 497.56668
$ws.Range("I80").Value = 532.63635
$ws.Range("J80").Value = 477.26315
$ws.Range("K80").Value = 532.63635
$ws.Range("L80").Value = 477.26315
$ws.Range("M80").Value = 465.36365
$ws.Range("N80").Value = -2473.26315

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 497.56668
$ws.Range("I83").Value = 532.63635
$ws.Range("J83").Value = 477.26315
$ws.Range("K83").Value = 2663.18175
$ws.Range("L83").Value = 2386.31575
$ws.Range("M83").Value = 2328.81825
$ws.Range("N83").Value = -12370.31575

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4001494
$ws.Range("I86").Value = 6251565
$ws.Range("K86").Value = 6251565
$ws.Range("M86").Value = -6250442

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4001494
$ws.Range("I89").Value = 6251565
$ws.Range("K89").Value = 31257825
$ws.Range("M89").Value = -31252209

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4331882.5
$ws.Range("I99").Value = 5293423
$ws.Range("K99").Value = 5293423
$ws.Range("M99").Value = -5291925

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3054.625
$ws.Range("I134").Value = 1388.4062
$ws.Range("K134").Value = 4165.2186
$ws.Range("M134").Value = -1630.2186

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23166.318
$ws.Range("J31").Value = 42969.668
$ws.Range("L31").Value = 42969.668
$ws.Range("N31").Value = -43559.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 23166.318
$ws.Range("J34").Value = 42969.668
$ws.Range("L34").Value = 42969.668
$ws.Range("N34").Value = -43373.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 65005.387
$ws.Range("I132").Value = 40714.035
$ws.Range("K132").Value = 122142.105
$ws.Range("M132").Value = -119612.105

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 5638
$ws.Range("I81").Value = 524
$ws.Range("J81").Value = 6774.4443
$ws.Range("K81").Value = 1572
$ws.Range("L81").Value = 20323.3329
$ws.Range("M81").Value = -449
$ws.Range("N81").Value = -22569.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 5638
$ws.Range("I84").Value = 524
$ws.Range("J84").Value = 6774.4443
$ws.Range("K84").Value = 4716
$ws.Range("L84").Value = 60969.9987
$ws.Range("M84").Value = 900
$ws.Range("N84").Value = -72201.9987

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1196.3334
$ws.Range("J107").Value = 1687
$ws.Range("L107").Value = 5061
$ws.Range("N107").Value = -8901

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 21934156
$ws.Range("J131").Value = 20837308
$ws.Range("L131").Value = 62511924
$ws.Range("N131").Value = -62522004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").ClearContents()
$ws.Range("N44").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 43916.8
$ws.Range("J95").Value = 43916.8
$ws.Range("L95").Value = 43916.8
$ws.Range("N95").Value = -49408.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 66166.664
$ws.Range("J108").Value = 71000
$ws.Range("L108").Value = 71000
$ws.Range("N108").Value = -78680

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4507056
$ws.Range("I113").Value = 7937996.5
$ws.Range("J113").Value = 3946.5
$ws.Range("K113").Value = 7937996.5
$ws.Range("L113").Value = 3946.5
$ws.Range("M113").Value = -7935826.5
$ws.Range("N113").Value = -8286.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3406.6206
$ws.Range("I132").Value = 3035.1304
$ws.Range("J132").Value = 4830.6665
$ws.Range("K132").Value = 9105.3912
$ws.Range("L132").Value = 14491.9995
$ws.Range("M132").Value = -6575.3912
$ws.Range("N132").Value = -19551.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 82155.09
$ws.Range("I22").Value = 444856.5
$ws.Range("J22").Value = 1554.7778
$ws.Range("K22").Value = 444856.5
$ws.Range("L22").Value = 1554.7778
$ws.Range("M22").Value = -444561.5
$ws.Range("N22").Value = -2144.7778

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 82155.09
$ws.Range("I27").Value = 444856.5
$ws.Range("J27").Value = 1554.7778
$ws.Range("K27").Value = 444856.5
$ws.Range("L27").Value = 1554.7778
$ws.Range("M27").Value = -444749.5
$ws.Range("N27").Value = -1768.7778

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 15876208
$ws.Range("I61").Value = 22223204
$ws.Range("K61").Value = 22223204
$ws.Range("M61").Value = -22223002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 15876208
$ws.Range("I113").Value = 22223204
$ws.Range("K113").Value = 22223204
$ws.Range("M113").Value = -22221034

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 30927.46
$ws.Range("I136").Value = 51252.43
$ws.Range("J136").Value = 7215
$ws.Range("K136").Value = 153757.29
$ws.Range("L136").Value = 21645
$ws.Range("M136").Value = -151207.29
$ws.Range("N136").Value = -26745

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4376.231
$ws.Range("I96").Value = 4262.909
$ws.Range("J96").Value = 4999.5
$ws.Range("K96").Value = 4262.909
$ws.Range("L96").Value = 4999.5
$ws.Range("M96").Value = -2889.909
$ws.Range("N96").Value = -7745.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 41667668
$ws.Range("I107").Value = 58823870
$ws.Range("J107").Value = 2606.1428
$ws.Range("K107").Value = 176471610
$ws.Range("L107").Value = 7818.428400000001
$ws.Range("M107").Value = -176469690
$ws.Range("N107").Value = -11658.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2270.5881
$ws.Range("I126").Value = 2512
$ws.Range("K126").Value = 7536
$ws.Range("M126").Value = -5066
